$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume (E) columns keep their text formatting so
# numeric-looking strings (e.g. "0.3620", "20.554.91") are not coerced to numbers.
$textCells = @(
    "D2",
    "E2",
    "D3",
    "E3",
    "E4",
    "D5",
    "E5",
    "D6",
    "E6",
    "D7",
    "E7",
    "D8",
    "E8",
    "D9",
    "E9",
    "E10",
    "D11",
    "E11",
    "E12",
    "D13",
    "E13",
    "D14",
    "E14",
    "D15",
    "E15",
    "D16",
    "E16",
    "D17",
    "E17",
    "D18",
    "E18",
    "D19",
    "E19",
    "D20",
    "E20",
    "D21",
    "E21",
    "E22",
    "D23",
    "E23",
    "D24",
    "E24",
    "D25",
    "E25",
    "D26",
    "E26",
    "D27",
    "E27",
    "D28",
    "E28",
    "D29",
    "E29",
    "D30",
    "E30",
    "D31",
    "E31",
    "D32",
    "E32",
    "D33",
    "E33",
    "D34",
    "E34",
    "D35",
    "E35",
    "D36",
    "E36",
    "D37",
    "E37",
    "D38",
    "E38",
    "D39",
    "E39",
    "D40",
    "E40",
    "D41",
    "E41",
    "D43",
    "E43",
    "D44",
    "E44",
    "D45",
    "E45",
    "D46",
    "E46",
    "D47",
    "E47",
    "D48",
    "E48",
    "D49",
    "E49",
    "D50",
    "E50",
    "D51",
    "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "20.554.91"
$ws.Range("E2").Value = "  +1.96%  "
$ws.Range("D3").Value = "1.474.53"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "0.9572"
$ws.Range("E5").Value = "  +4.99%  "
$ws.Range("D6").Value = "277.95"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "0.3620"
$ws.Range("E7").Value = "  -1.10%  "
$ws.Range("D8").Value = "0.3070"
$ws.Range("E8").Value = "  -1.82%  "
$ws.Range("D9").Value = "39.67"
$ws.Range("E9").Value = "  +2.41%  "
$ws.Range("E10").Value = "  +4.64%  "
$ws.Range("D11").Value = "0.06658"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "5.539"
$ws.Range("E13").Value = "  +2.75%  "
$ws.Range("D14").Value = "18.13"
$ws.Range("E14").Value = "  +3.55%  "
$ws.Range("D15").Value = "6.196"
$ws.Range("E15").Value = "  +2.12%  "
$ws.Range("D16").Value = "0.9569"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("D17").Value = "0.00001027"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").Value = "1.475.72"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").Value = "0.05926"
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("D20").Value = "69.22"
$ws.Range("E20").Value = "  +2.24%  "
$ws.Range("D21").Value = "5.506"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "11.19"
$ws.Range("E23").Value = "  +3.58%  "
$ws.Range("D24").Value = "2.256"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "20.568.32"
$ws.Range("E25").Value = "  +1.91%  "
$ws.Range("D26").Value = "143.53"
$ws.Range("E26").Value = "  +5.28%  "
$ws.Range("D27").Value = "2.129"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("D28").Value = "17.19"
$ws.Range("E28").Value = "  +1.21%  "
$ws.Range("D29").Value = "1.638.22"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("D30").Value = "113.76"
$ws.Range("E30").Value = "  +2.65%  "
$ws.Range("D31").Value = "3.920"
$ws.Range("E31").Value = "  +4.92%  "
$ws.Range("D32").Value = "4.991"
$ws.Range("E32").Value = "  +3.58%  "
$ws.Range("D33").Value = "0.8094"
$ws.Range("E33").Value = "  +0.42%  "
$ws.Range("D34").Value = "0.07998"
$ws.Range("E34").Value = "  +3.88%  "
$ws.Range("D35").Value = "1.516"
$ws.Range("E35").Value = "  +4.18%  "
$ws.Range("D36").Value = "1.219"
$ws.Range("E36").Value = "  +7.48%  "
$ws.Range("D37").Value = "0.05787"
$ws.Range("E37").Value = "  -3.49%  "
$ws.Range("D38").Value = "4.745"
$ws.Range("E38").Value = "  +1.14%  "
$ws.Range("D39").Value = "0.02058"
$ws.Range("E39").Value = "  +3.29%  "
$ws.Range("D40").Value = "10.40"
$ws.Range("E40").Value = "  +2.51%  "
$ws.Range("D41").Value = "0.9576"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("D43").Value = "7.435"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("D44").Value = "0.5286"
$ws.Range("E44").Value = "  +1.04%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.29"
$ws.Range("E45").Value = "  +2.47%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").Value = "3.523"
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").Value = "118.35"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("D48").Value = "0.5207"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("D49").Value = "1.818"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "0.06480"
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("D51").Value = "0.9852"
$ws.Range("E51").Value = "  -1.12%  "
